$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 - MagneticSensor: claim it (add new-name, developer, finished, mode, interface)
$ws.Range("F44").Value2 = "Magnetic"
$ws.Range("B44").Value2 = "HiTechnicMagneticSensor"
$ws.Range("D44").Value2 = "Lawrie"
$ws.Range("E44").Value2 = "N"
$ws.Range("G44").Value2 = "SampleProvider"

# Row 61 - RCXLightSensor: claim it
$ws.Range("B61").Value2 = "RCXLightSensor"
$ws.Range("D61").Value2 = "Lawrie"
$ws.Range("E61").Value2 = "N"
$ws.Range("F61").Value2 = "Light"
$ws.Range("G61").Value2 = "SampleProvider"

# Rows 62-68 - mark "Fits in framework" column
$ws.Range("C62").Value2 = "N"
$ws.Range("C63").Value2 = "N"
$ws.Range("C64").Value2 = "N"
$ws.Range("C65").Value2 = "N"
$ws.Range("C66").Value2 = "N"
$ws.Range("C67").Value2 = "?"
$ws.Range("C68").Value2 = "N"

# Update active cell selection to C68
$ws.Range("C68").Select()
